$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Grupo A')
$ws.Cells.Item(2, 2).Value = 'Dom Camillo68'
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 520.8798828125
$ws.Cells.Item(2, 8).Value = 467.770263671875
$ws.Cells.Item(2, 9).Value = 53.109619140625
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'pura bucha /botafogo'
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 466.9794921875
$ws.Cells.Item(3, 8).Value = 441.240234375
$ws.Cells.Item(3, 9).Value = 25.7392578125
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'BORGES ITAQUI F.C.'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 475.830078125
$ws.Cells.Item(4, 8).Value = 509.7998046875
$ws.Cells.Item(4, 9).Value = -33.9697265625
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'A Lenda Super Vasco F.c'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(5, 7).Value = 472.790283203125
$ws.Cells.Item(5, 8).Value = 517.66943359375
$ws.Cells.Item(5, 9).Value = -44.879150390625
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo B')
$ws.Cells.Item(2, 2).Value = 'Tabajara de Inhaua FC2'
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 522.12060546875
$ws.Cells.Item(2, 8).Value = 480.06005859375
$ws.Cells.Item(2, 9).Value = 42.060546875
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'lsauer fc'
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 513.55029296875
$ws.Cells.Item(3, 8).Value = 487.76025390625
$ws.Cells.Item(3, 9).Value = 25.7900390625
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'HS SPORTS F.C'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 498.68017578125
$ws.Cells.Item(4, 8).Value = 498.5908203125
$ws.Cells.Item(4, 9).Value = 0.08935546875
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'Rolo Compressor ZN'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(5, 7).Value = 410.4404296875
$ws.Cells.Item(5, 8).Value = 478.38037109375
$ws.Cells.Item(5, 9).Value = -67.93994140625
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo C')
$ws.Cells.Item(2, 2).Value = 'Analove10 ITAQUI GRANDE!!'
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 526.31005859375
$ws.Cells.Item(2, 8).Value = 499.45947265625
$ws.Cells.Item(2, 9).Value = 26.8505859375
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'TEAM LOPES 99'
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 518.08984375
$ws.Cells.Item(3, 8).Value = 501.429931640625
$ws.Cells.Item(3, 9).Value = 16.659912109375
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'Grêmio imortal 37'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 517.379638671875
$ws.Cells.Item(4, 8).Value = 497.919921875
$ws.Cells.Item(4, 9).Value = 19.459716796875
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'cartola scheuer'
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = 491.15966796875
$ws.Cells.Item(5, 8).Value = 554.1298828125
$ws.Cells.Item(5, 9).Value = -62.97021484375
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo D')
$ws.Cells.Item(2, 2).Value = 'Texas Club 2025'
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 501.6298828125
$ws.Cells.Item(2, 8).Value = 477.0302734375
$ws.Cells.Item(2, 9).Value = 24.599609375
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'Super Vasco f.c'
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 505.749755859375
$ws.Cells.Item(3, 8).Value = 486.06982421875
$ws.Cells.Item(3, 9).Value = 19.679931640625
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'Tatols Beants F.C'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 503.4296875
$ws.Cells.Item(4, 8).Value = 531.60009765625
$ws.Cells.Item(4, 9).Value = -28.17041015625
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'Fedato Futebol Clube'
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = 493.73046875
$ws.Cells.Item(5, 8).Value = 509.839599609375
$ws.Cells.Item(5, 9).Value = -16.109130859375
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo E')
$ws.Cells.Item(2, 2).Value = 'Real SCI'
$ws.Cells.Item(2, 3).Value = 15
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 552.1396484375
$ws.Cells.Item(2, 8).Value = 445.080322265625
$ws.Cells.Item(2, 9).Value = 107.059326171875
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'Gremiomaniasm'
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 500.619384765625
$ws.Cells.Item(3, 8).Value = 454.81982421875
$ws.Cells.Item(3, 9).Value = 45.799560546875
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'E.C. Bororé'
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(4, 7).Value = 479.160400390625
$ws.Cells.Item(4, 8).Value = 537.1396484375
$ws.Cells.Item(4, 9).Value = -57.979248046875
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'ITAQUI F. C.'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(5, 7).Value = 398.92041015625
$ws.Cells.Item(5, 8).Value = 493.800048828125
$ws.Cells.Item(5, 9).Value = -94.879638671875
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo F')
$ws.Cells.Item(2, 2).Value = 'TORRESMO COM PINGA'
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 514.140625
$ws.Cells.Item(2, 8).Value = 447.199951171875
$ws.Cells.Item(2, 9).Value = 66.940673828125
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'Lá do Itaqui'
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 507.43017578125
$ws.Cells.Item(3, 8).Value = 474.50048828125
$ws.Cells.Item(3, 9).Value = 32.9296875
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'FC Los Castilho'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 476.580078125
$ws.Cells.Item(4, 8).Value = 498.669921875
$ws.Cells.Item(4, 9).Value = -22.08984375
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'seralex'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(5, 7).Value = 393.249755859375
$ws.Cells.Item(5, 8).Value = 471.0302734375
$ws.Cells.Item(5, 9).Value = -77.780517578125
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo G')
$ws.Cells.Item(2, 2).Value = 'KING LEONN'
$ws.Cells.Item(2, 3).Value = 18
$ws.Cells.Item(2, 4).Value = 6
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 513.68994140625
$ws.Cells.Item(2, 8).Value = 411.539306640625
$ws.Cells.Item(2, 9).Value = 102.150634765625
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'KillerColorado'
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 469.429931640625
$ws.Cells.Item(3, 8).Value = 467.539794921875
$ws.Cells.Item(3, 9).Value = 1.89013671875
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'F.C. Rei Das Copas'
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(4, 7).Value = 449.68017578125
$ws.Cells.Item(4, 8).Value = 449.2294921875
$ws.Cells.Item(4, 9).Value = 0.45068359375
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'TATITTA FC'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(5, 7).Value = 408.619384765625
$ws.Cells.Item(5, 8).Value = 513.11083984375
$ws.Cells.Item(5, 9).Value = -104.491455078125
$ws.Cells.Item(5, 10).Value = 4

$ws = $wb.Worksheets.Item('Grupo H')
$ws.Cells.Item(2, 2).Value = 'Gig@ntte'
$ws.Cells.Item(2, 3).Value = 12
$ws.Cells.Item(2, 4).Value = 4
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 504.9794921875
$ws.Cells.Item(2, 8).Value = 485.520263671875
$ws.Cells.Item(2, 9).Value = 19.459228515625
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(3, 2).Value = 'Laranjja Mecannica'
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 443.39990234375
$ws.Cells.Item(3, 8).Value = 466.25
$ws.Cells.Item(3, 9).Value = -22.85009765625
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(4, 2).Value = 'MauHumor F.C.'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 540.13037109375
$ws.Cells.Item(4, 8).Value = 485.1298828125
$ws.Cells.Item(4, 9).Value = 55.00048828125
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(5, 2).Value = 'FBC Colorado'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(5, 7).Value = 458.259765625
$ws.Cells.Item(5, 8).Value = 509.869384765625
$ws.Cells.Item(5, 9).Value = -51.609619140625
$ws.Cells.Item(5, 10).Value = 4

